$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("B2").Value = 255.895695
$ws.Range("D2").Value = 1033.577748

$ws.Range("B3").Value = 4.31741
$ws.Range("D3").Value = 8.719137
$ws.Range("E3").Value = 0.000226

$ws.Range("B4").Value = 55.210883
$ws.Range("C4").Value = 223

$ws.Range("G5").Value = 0.031179
$ws.Range("H5").Value = -0.186131
$ws.Range("I5").Value = 0.248489
$ws.Range("J5").Value = 0.938802

$ws.Range("G6").Value = -0.273449
$ws.Range("H6").Value = -0.507446
$ws.Range("I6").Value = -0.039453
$ws.Range("J6").Value = 0.017291

$ws.Range("G7").Value = -0.304628
$ws.Range("H7").Value = -0.481532
$ws.Range("I7").Value = -0.127724
$ws.Range("J7").Value = 0.000198
